$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2841271
$ws.Range("I12").Value = 3788003.2
$ws.Range("K12").Value = 3788003.2
$ws.Range("M12").Value = -3787833.2

$ws.Range("H39").Value = 416.55554
$ws.Range("I39").Value = 1670
$ws.Range("J39").Value = 259.875
$ws.Range("K39").Value = 5010
$ws.Range("L39").Value = 779.625
$ws.Range("M39").Value = -4714
$ws.Range("N39").Value = -1371.625

$ws.Range("H70").Value = 7166.1665
$ws.Range("J70").Value = 8399.4
$ws.Range("L70").Value = 25198.2
$ws.Range("N70").Value = -25738.2

$ws.Range("H73").Value = 7166.1665
$ws.Range("J73").Value = 8399.4
$ws.Range("L73").Value = 25198.2
$ws.Range("N73").Value = -27070.2

$ws.Range("H106").Value = 100006000
$ws.Range("J106").Value = 3500
$ws.Range("L106").Value = 3500
$ws.Range("N106").Value = -4762

$ws.Range("H118").Value = 100000616
$ws.Range("I118").Value = 142857540
$ws.Range("K118").Value = 428572620
$ws.Range("M118").Value = -428570963

$ws.Range("H132").Value = 2836.9487
$ws.Range("I132").Value = 3017.5588
$ws.Range("K132").Value = 9052.6764
$ws.Range("M132").Value = -6522.6764

$ws.Range("H137").Value = 34765.742
$ws.Range("I137").Value = 52083.973
$ws.Range("J137").Value = 2863.7368
$ws.Range("K137").Value = 156251.919
$ws.Range("L137").Value = 8591.2104
$ws.Range("M137").Value = -153701.919
$ws.Range("N137").Value = -13691.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1010790.25
$ws.Range("I2").Value = 1131887.2
$ws.Range("K2").Value = 1131887.2
$ws.Range("M2").Value = -1131774.2

$ws.Range("H32").Value = 10763.564
$ws.Range("I32").Value = 7054.3623
$ws.Range("J32").Value = 18731.482
$ws.Range("K32").Value = 7054.3623
$ws.Range("L32").Value = 18731.482
$ws.Range("M32").Value = -6767.3623
$ws.Range("N32").Value = -19305.482

$ws.Range("H45").Value = 6806119
$ws.Range("I45").Value = 10990640
$ws.Range("J45").Value = 6273.5
$ws.Range("K45").Value = 10990640
$ws.Range("L45").Value = 6273.5
$ws.Range("M45").Value = -10990263
$ws.Range("N45").Value = -7027.5

$ws.Range("H47").Value = 39999.668
$ws.Range("J47").Value = 39999.668
$ws.Range("L47").Value = 39999.668
$ws.Range("N47").Value = -41449.668

$ws.Range("H61").Value = 3755.4333
$ws.Range("I61").Value = 3780.3635
$ws.Range("K61").Value = 3780.3635
$ws.Range("M61").Value = -3568.3635

$ws.Range("H74").Value = 22674.773
$ws.Range("J74").Value = 103429.89
$ws.Range("L74").Value = 103429.89
$ws.Range("N74").Value = -105177.89

$ws.Range("H77").Value = 22674.773
$ws.Range("J77").Value = 103429.89
$ws.Range("L77").Value = 517149.45
$ws.Range("N77").Value = -525885.45

$ws.Range("H116").Value = 1010790.25
$ws.Range("I116").Value = 1131887.2
$ws.Range("K116").Value = 1131887.2
$ws.Range("M116").Value = -1129593.2

$ws.Range("H122").Value = 1599555.9
$ws.Range("I122").Value = 1755689.2
$ws.Range("K122").Value = 5267067.6
$ws.Range("M122").Value = -5264617.6

$ws.Range("H132").Value = 1980.0476
$ws.Range("I132").Value = 1050.8064
$ws.Range("K132").Value = 3152.4192
$ws.Range("M132").Value = -622.4191999999998

$ws.Range("H136").Value = 3755.4333
$ws.Range("I136").Value = 3780.3635
$ws.Range("K136").Value = 11341.0905
$ws.Range("M136").Value = -8791.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1058013.4
$ws.Range("I3").Value = 1179041
$ws.Range("J3").Value = 1772.2727
$ws.Range("K3").Value = 1179041
$ws.Range("L3").Value = 1772.2727
$ws.Range("M3").Value = -1178927
$ws.Range("N3").Value = -2000.2727

$ws.Range("H94").Value = 6257354.5
$ws.Range("I94").Value = 11113965
$ws.Range("K94").Value = 11113965
$ws.Range("M94").Value = -11113514

$ws.Range("H105").Value = 5684137.5
$ws.Range("I105").Value = 5684137.5
$ws.Range("K105").Value = 5684137.5
$ws.Range("M105").Value = -5682390.5

$ws.Range("H107").Value = 8932205
$ws.Range("I107").Value = 11907120
$ws.Range("K107").Value = 11907120
$ws.Range("M107").Value = -11905200

$ws.Range("H132").Value = 83799.39999999999
$ws.Range("J132").Value = 83799.39999999999
$ws.Range("L132").Value = 83799.39999999999
$ws.Range("N132").Value = -93919.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1376.6
$ws.Range("I105").Value = 1327.6666
$ws.Range("J105").Value = 1450
$ws.Range("K105").Value = 1327.6666
$ws.Range("L105").Value = 1450
$ws.Range("M105").Value = 419.3334
$ws.Range("N105").Value = -4944

$ws.Range("H107").Value = 1396.2766
$ws.Range("I107").Value = 1421.6216
$ws.Range("J107").Value = 1302.5
$ws.Range("K107").Value = 1421.6216
$ws.Range("L107").Value = 1302.5
$ws.Range("M107").Value = 498.3784000000001
$ws.Range("N107").Value = -5142.5

$ws.Range("H134").Value = 1684.8918
$ws.Range("I134").Value = 1101.3939
$ws.Range("K134").Value = 3304.1817
$ws.Range("M134").Value = -769.1817000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2745.2666
$ws.Range("J132").Value = 2998.2727
$ws.Range("L132").Value = 26984.4543
$ws.Range("N132").Value = -32044.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 36757.855
$ws.Range("I45").Value = 27666.666
$ws.Range("J45").Value = 53122
$ws.Range("K45").Value = 27666.666
$ws.Range("L45").Value = 53122
$ws.Range("M45").Value = -27107.666
$ws.Range("N45").Value = -54240

$ws.Range("H80").Value = 24990364
$ws.Range("I80").Value = 47688924
$ws.Range("J80").Value = 21945.2
$ws.Range("K80").Value = 47688924
$ws.Range("L80").Value = 21945.2
$ws.Range("M80").Value = -47687926
$ws.Range("N80").Value = -23941.2

$ws.Range("H83").Value = 24990364
$ws.Range("I83").Value = 47688924
$ws.Range("J83").Value = 21945.2
$ws.Range("K83").Value = 238444620
$ws.Range("L83").Value = 109726
$ws.Range("M83").Value = -238439628
$ws.Range("N83").Value = -119710

$ws.Range("H97").Value = 794315.25
$ws.Range("I97").Value = 1035746.56
$ws.Range("J97").Value = 1040.8572
$ws.Range("K97").Value = 1035746.56
$ws.Range("L97").Value = 1040.8572
$ws.Range("M97").Value = -1035250.56
$ws.Range("N97").Value = -2032.8572

$ws.Range("H113").Value = 5558128.5
$ws.Range("I113").Value = 7248255
$ws.Range("K113").Value = 7248255
$ws.Range("M113").Value = -7246085

$ws.Range("H122").Value = 408180.53
$ws.Range("I122").Value = 812169.4399999999
$ws.Range("J122").Value = 4191.636
$ws.Range("K122").Value = 2436508.32
$ws.Range("L122").Value = 12574.908
$ws.Range("M122").Value = -2434058.32
$ws.Range("N122").Value = -17474.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6184.095
$ws.Range("I40").Value = 3498.4285
$ws.Range("J40").Value = 11555.429
$ws.Range("K40").Value = 3498.4285
$ws.Range("L40").Value = 11555.429
$ws.Range("M40").Value = -3362.4285
$ws.Range("N40").Value = -11827.429

$ws.Range("H55").Value = 2463.7273
$ws.Range("I55").Value = 2522.6667
$ws.Range("J55").Value = 2441.625
$ws.Range("K55").Value = 2522.6667
$ws.Range("L55").Value = 2441.625
$ws.Range("M55").Value = -2349.6667
$ws.Range("N55").Value = -2787.625

$ws.Range("H61").Value = 3704369.2
$ws.Range("I61").Value = 4115791
$ws.Range("K61").Value = 4115791
$ws.Range("M61").Value = -4115589

$ws.Range("H82").Value = 50266056
$ws.Range("I82").Value = 65973850
$ws.Range("J82").Value = 1121.8
$ws.Range("K82").Value = 65973850
$ws.Range("L82").Value = 1121.8
$ws.Range("M82").Value = -65973489
$ws.Range("N82").Value = -1843.8

$ws.Range("H85").Value = 50266056
$ws.Range("I85").Value = 65973850
$ws.Range("J85").Value = 1121.8
$ws.Range("K85").Value = 65973850
$ws.Range("L85").Value = 1121.8
$ws.Range("M85").Value = -65972602
$ws.Range("N85").Value = -3617.8

$ws.Range("H93").Value = 16678494
$ws.Range("I93").Value = 22223820
$ws.Range("K93").Value = 22223820
$ws.Range("M93").Value = -22222572

$ws.Range("H113").Value = 3704369.2
$ws.Range("I113").Value = 4115791
$ws.Range("K113").Value = 4115791
$ws.Range("M113").Value = -4113621

$ws.Range("H132").Value = 6480.12
$ws.Range("I132").Value = 7106.263
$ws.Range("J132").Value = 4497.3335
$ws.Range("K132").Value = 21318.789
$ws.Range("L132").Value = 13492.0005
$ws.Range("M132").Value = -18788.789
$ws.Range("N132").Value = -18552.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H29").Value = 11999.6
$ws.Range("J29").Value = 14999
$ws.Range("L29").Value = 14999
$ws.Range("N29").Value = -15579

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H113").Value = 1175.8
$ws.Range("I113").Value = 267.93332
$ws.Range("J113").Value = 2537.6
$ws.Range("K113").Value = 803.7999599999999
$ws.Range("L113").Value = 7612.799999999999
$ws.Range("M113").Value = 1366.20004
$ws.Range("N113").Value = -11952.8

$ws.Range("H122").Value = 3950.077
$ws.Range("I122").Value = 2264.3157
$ws.Range("K122").Value = 6792.9471
$ws.Range("M122").Value = -4342.9471

$ws.Range("H136").Value = 884.34784
$ws.Range("I136").Value = 794.9322
$ws.Range("K136").Value = 2384.7966
$ws.Range("M136").Value = 165.2034000000003
